$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Row 4 - column C (reuses existing shared string)
$ws.Range("C4").Value = "Biber only looks for perfect aspect, we draw a finer distinction between present and past perfect. If we want to we can easily combine them later"

# Row 5 - column C ("same" - new shared string)
$ws.Range("C5").Value = "same"

# Row 28 - columns C and D
$ws.Range("C28").Value = "same"
$ws.Range("D28").Value = "appears to work well"

# Row 29 - column C
$ws.Range("C29").Value = "same"

# Row 30 - columns C and D
$ws.Range("C30").Value = "same"
$ws.Range("D30").Value = "not good, catches unintended stuff like ""with prices going up"" or ""or is the passage saying something quite different"" - might need to drop this one"

# Row 29 - column D (new shared string, added after row 30's D text)
$ws.Range("D29").Value = "currently can't check this due to tagger problems (tags the VBNs we are looking for as VBD)"

# Row 5 - column D (new shared string, added last)
$ws.Range("D5").Value = "some tagger inaccuracy"

# Row 60 - column D
$ws.Range("D60").Value = "works well"

# Rows 77-82 - column C
$ws.Range("C77").Value = "(not included in Biber)"
$ws.Range("C78").Value = "(not included in Biber)"
$ws.Range("C79").Value = "(not included in Biber)"
$ws.Range("C80").Value = "(not included in Biber)"
$ws.Range("C81").Value = "(not included in Biber)"
$ws.Range("C82").Value = "(not included in Biber)"

# Update the selection to match the latest view state
$ws.Activate()
$ws.Range("D61").Select()
